$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算  (A1:F16 -> A1:F17)
# New row 2: 日期：2021/12/08 / 202201 / 17793 / 11174 / 3807702 / 17647
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("台指期換倉成本計算")
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2").Value = "日期：2021/12/08"
$ws1.Range("B3").Copy($ws1.Range("B2"))
$ws1.Range("C2").Value = 17793
$ws1.Range("D2").Value = 11174
$ws1.Range("E2").Value = 3807702
$ws1.Range("F2").Value = 17647

# ---------------------------------------------------------------
# Sheet 2: 散戶多空力道  (A1:B31 -> A1:B32)
# New row 2: 日期：2021/12/08 / -0.03
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("散戶多空力道")
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = "日期：2021/12/08"
$ws2.Range("B2").Value = -0.03

# ---------------------------------------------------------------
# Sheet 3: 三大法人買賣金額  (A1:C31 -> A1:C32)
# New row 2: 110年12月08日 / 169.47 / -133.97
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("三大法人買賣金額")
$ws3.Rows.Item(2).Insert()
$ws3.Range("A2").Value = "110年12月08日"
$ws3.Range("B2").Value = 169.47
$ws3.Range("C2").Value = -133.97

# ---------------------------------------------------------------
# Sheet 4: 大盤多空點位  (A1:B30 -> A1:B31)
# New row 2: 110年12月08日 / 17913.87
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("大盤多空點位")
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = "110年12月08日"
$ws4.Range("B2").Value = 17913.87

# ---------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位  (A1:N29 -> A1:N30)
# New row 2: 2021/12/08 / 47996 / 55473 / -2741 / -3078 / 25001 / 50415 /
#            -349 / -1153 / -25414 / 804 / -2392 / -1925 / -467
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("期貨大額交易人未沖銷部位")
$ws5.Rows.Item(2).Insert()
$ws5.Range("A2").Value = "'2021/12/08"
$ws5.Range("A3").Copy()
$ws5.Range("A2").PasteSpecial(-4122)
$ws5.Range("B2").Value = 47996
$ws5.Range("C2").Value = 55473
$ws5.Range("D2").Value = -2741
$ws5.Range("E2").Value = -3078
$ws5.Range("F2").Value = 25001
$ws5.Range("G2").Value = 50415
$ws5.Range("H2").Value = -349
$ws5.Range("I2").Value = -1153
$ws5.Range("J2").Value = -25414
$ws5.Range("K2").Value = 804
$ws5.Range("L2").Value = -2392
$ws5.Range("M2").Value = -1925
$ws5.Range("N2").Value = -467
